$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldChemo = 'tp.chemotherapy_regimen In ["Other treatment given as part of a CTSU protocol"] and d.er_status'
$newChemo = 'tp.chemotherapy_regimen In ["Dose dense AC (2 week cycles)"]  and d.er_status'

# B2: Cases query (long "demographic" query), B4: Files query, C2/C3/C4: Stats query
$targets = @("B2", "C2", "B4", "C3", "C4")

foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    $val = $cell.Value2
    if ($val -ne $null -and $val.Contains($oldChemo)) {
        $cell.Value2 = $val.Replace($oldChemo, $newChemo)
    }
}

# Update the worksheet view: the saved view now has the selection on C4
# (matches the saved view state in the target file).
$ws.Range("C4").Select()
